$d = $word.ActiveDocument

# This revision appends a reviewer's name to the "Reviewer(s):" line in
# the document's header block, turning:
#   "Reviewer(s):"
# into:
#   "Reviewer(s): (Lukas)"
#
# Locate that paragraph explicitly by its text (rather than assuming
# it is simply the last paragraph) so the edit is resilient to the
# exact layout of the rest of the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Reviewer(s):*") {
        $target = $p
    }
}

if ($target -eq $null) {
    $target = $d.Paragraphs.Last
}

$r = $target.Range
$r.Collapse(0)
$r.InsertAfter(" (Lukas)")
